$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-06 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-07 Saturday", 2)

$d.Content.Find.Execute("956÷2=478, 0", $true, $false, $false, $false, $false, $true, 1, $false, "160÷4=40, 0", 2)
$d.Content.Find.Execute("659÷6=109, 5", $true, $false, $false, $false, $false, $true, 1, $false, "642÷4=160, 2", 2)
$d.Content.Find.Execute("758÷5=151, 3", $true, $false, $false, $false, $false, $true, 1, $false, "657÷7=93, 6", 2)
$d.Content.Find.Execute("555÷3=185, 0", $true, $false, $false, $false, $false, $true, 1, $false, "569÷7=81, 2", 2)
$d.Content.Find.Execute("783÷4=195, 3", $true, $false, $false, $false, $false, $true, 1, $false, "155÷6=25, 5", 2)

$d.Content.Find.Execute("860÷5=172, 0", $true, $false, $false, $false, $false, $true, 1, $false, "778÷9=86, 4", 2)
$d.Content.Find.Execute("108÷4=27, 0", $true, $false, $false, $false, $false, $true, 1, $false, "654÷6=109, 0", 2)
$d.Content.Find.Execute("350÷4=87, 2", $true, $false, $false, $false, $false, $true, 1, $false, "858÷5=171, 3", 2)
$d.Content.Find.Execute("724÷3=241, 1", $true, $false, $false, $false, $false, $true, 1, $false, "904÷2=452, 0", 2)
$d.Content.Find.Execute("749÷6=124, 5", $true, $false, $false, $false, $false, $true, 1, $false, "743÷6=123, 5", 2)

$d.Content.Find.Execute("657÷8=82, 1", $true, $false, $false, $false, $false, $true, 1, $false, "777÷9=86, 3", 2)
$d.Content.Find.Execute("669÷9=74, 3", $true, $false, $false, $false, $false, $true, 1, $false, "580÷4=145, 0", 2)
$d.Content.Find.Execute("327÷5=65, 2", $true, $false, $false, $false, $false, $true, 1, $false, "915÷2=457, 1", 2)
$d.Content.Find.Execute("969÷7=138, 3", $true, $false, $false, $false, $false, $true, 1, $false, "662÷4=165, 2", 2)
$d.Content.Find.Execute("437÷8=54, 5", $true, $false, $false, $false, $false, $true, 1, $false, "483÷3=161, 0", 2)

$d.Content.Find.Execute("146÷4=36, 2", $true, $false, $false, $false, $false, $true, 1, $false, "443÷2=221, 1", 2)
$d.Content.Find.Execute("183÷9=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "398÷6=66, 2", 2)
$d.Content.Find.Execute("219÷8=27, 3", $true, $false, $false, $false, $false, $true, 1, $false, "261÷7=37, 2", 2)
$d.Content.Find.Execute("222÷9=24, 6", $true, $false, $false, $false, $false, $true, 1, $false, "606÷4=151, 2", 2)
$d.Content.Find.Execute("170÷8=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "285÷8=35, 5", 2)

$d.Content.Find.Execute("343÷2=171, 1", $true, $false, $false, $false, $false, $true, 1, $false, "310÷6=51, 4", 2)
$d.Content.Find.Execute("429÷9=47, 6", $true, $false, $false, $false, $false, $true, 1, $false, "491÷7=70, 1", 2)
$d.Content.Find.Execute("236÷7=33, 5", $true, $false, $false, $false, $false, $true, 1, $false, "922÷6=153, 4", 2)
$d.Content.Find.Execute("842÷6=140, 2", $true, $false, $false, $false, $false, $true, 1, $false, "529÷4=132, 1", 2)
$d.Content.Find.Execute("132÷9=14, 6", $true, $false, $false, $false, $false, $true, 1, $false, "469÷2=234, 1", 2)

Write-Output "Replacements complete"
